$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 3710237
$ws.Range("C2").Value = 6630298
$ws.Range("C3").Value = 6401285
$ws.Range("C4").Value = 6612387

$ws.Range("D6").Select() | Out-Null
